$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
"48+11=59",
"21+30=51",
"26+63=89",
"81-20=61",
"15+67=82",
"81-54=27",
"72-72=0",
"25+16=41",
"26+16=42",
"48+37=85",
"77-9=68",
"33+48=81",
"29+7=36",
"68-20=48",
"42-38=4",
"67-42=25",
"69-64=5",
"66+22=88",
"68-24=44",
"76-46=30",
"60-22=38",
"27+24=51",
"33+50=83",
"36-5=31",
"92-81=11",
"11+84=95",
"86-57=29",
"1+9=10",
"44+45=89",
"37-27=10",
"52-27=25",
"47+43=90",
"35-32=3",
"48-30=18",
"3+49=52",
"9+41=50",
"86-76=10",
"10+78=88",
"26+68=94",
"84-7=77",
"17-11=6",
"58-42=16",
"15+39=54",
"83-10=73",
"22+9=31",
"18+42=60",
"98-17=81",
"47+17=64",
"8+1=9",
"73-31=42",
"18+38=56",
"90-66=24",
"47+39=86",
"98-53=45",
"83-77=6",
"7+58=65",
"80-42=38",
"59+13=72",
"23+69=92",
"60-22=38",
"49-27=22",
"11+40=51",
"43-4=39",
"30+2=32",
"39+57=96",
"20+28=48",
"60+12=72",
"15-9=6",
"75-22=53",
"9+78=87",
"94-83=11",
"40-6=34",
"71-68=3",
"8+72=80",
"26-24=2",
"36+58=94",
"43+30=73",
"95-71=24",
"48+31=79",
"94-26=68",
"29+70=99",
"24+8=32",
"79-55=24",
"0+56=56",
"30+36=66",
"1+58=59",
"34+13=47",
"73+16=89",
"60-56=4",
"71+10=81",
"48-11=37",
"7+68=75",
"5+37=42",
"93-69=24",
"94+0=94",
"22+7=29",
"70+25=95",
"43+41=84",
"56-52=4",
"90-51=39"
)

$cols = $t.Columns.Count
$rows = $t.Rows.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
  for ($c = 1; $c -le $cols; $c++) {
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $values[$idx]
    $idx = $idx + 1
  }
}
Write-Output ("updated cells: " + $idx)
